# "Generate Report for Handback" - refresh the localization-status report
# after a successful handback: status flips from "Ready for handoff" to
# "Handed back: in sync with en-US", the handback timestamps advance, and
# the stale "handback file is not the latest" error clears now that the
# handback is in sync.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: per-language status columns -------------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2").Value = $newStatus
$ovw.Range("F2").Value = $newStatus

# --- zh-cn sheet -----------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-21 16:56:19"
$zhcn.Range("P2").Value = ""

# --- de-de sheet -----------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-21 16:56:25"
$dede.Range("P2").Value = ""

# --- column widths: widen the Status columns to fit the longer text, and
#     shrink the now-empty Error Detail columns. (ColumnWidth is stored on
#     a quantized character grid, so these land on the nearest achievable
#     grid value to the authored widths.)
$ovw.Columns.Item(5).ColumnWidth = 29.166666666666668
$ovw.Columns.Item(6).ColumnWidth = 29.166666666666668

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333334

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(16).ColumnWidth = 12.833333333333334
